$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at row 10 (shifts LAMIFEN..syringe rows down by one),
#    then copy formatting from the row above (row 9) so the new row matches
#    the table's styling (borders, fonts, text number-format, merges).
$ws.Rows.Item(10).Insert()
$ws.Range("A9:Q9").Copy()
$ws.Range("A10:Q10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Re-create the per-cell merges for the new row 10 (matching the pattern
# used by every other item row in the table).
$ws.Range("A10:B10").Merge()
$ws.Range("C10:G10").Merge()
$ws.Range("H10:K10").Merge()
$ws.Range("L10:M10").Merge()
$ws.Range("N10:O10").Merge()

# 2. Populate the new row 10 with the FLAGYL item data (item #4).
#    Leading apostrophes force text storage for numeric-looking values so
#    they match the workbook's existing "number stored as text" pattern.
$ws.Range("A10").Value2 = 4
$ws.Range("C10").Value2 = "FLAGYL 125MG/5ML 100 ML SUSPENSION"
$ws.Range("H10").Value2 = "'5:0"
$ws.Range("L10").Value2 = "'1"
$ws.Range("N10").Value2 = "'26.00"
$ws.Range("P10").Value2 = "'26.0000"
$ws.Range("Q10").Value2 = "'1:0"

# 3. Renumber the items that were pushed down (LAMIFEN..syringes) so the
#    "م" / item-number column stays sequential.
$ws.Range("A11").Value2 = 5
$ws.Range("A12").Value2 = 6
$ws.Range("A13").Value2 = 7
$ws.Range("A14").Value2 = 8
$ws.Range("A15").Value2 = 9

# 4. The last item (سرنجات 3 سم) now sits on row 14; its transaction-count
#    column changed from 3:0 to 1:0.
$ws.Range("Q14").Value2 = "'1:0"

# 5. Update the grand-total (sale price column) to include the new item.
$ws.Range("P16").Value2 = 367.8

# 6. Update the generated timestamp text in the footer.
$ws.Range("A17").Value2 = "Saturday, 2 August, 2025 10:30 AM"
